# Add a new "Swiss" worksheet (Switzerland market test data), modelled on
# the existing "Czech" sheet, and update the workbook/sheet selection state
# to match.

$wb = $excel.ActiveWorkbook
$czech = $wb.Worksheets.Item("Czech")

# Before losing focus, the Czech sheet's selection becomes a "select all"
# (A1:XFD1048576), matching the target selection state.
$czech.Range("A1:XFD1048576").Select()

# Duplicate the Czech sheet; Excel places the copy immediately after the
# source sheet and automatically activates it (making it the new
# tabSelected sheet / activeTab).
$czech.Copy($null, $czech)

# The copy becomes the new last sheet - rename it and update its market
# specific values.
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2343/T2641"
$swiss.Range("B5").Select()
